$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all Timestamp values (column A) forward by 22 days (add 22 to the date serial),
# and update the Notified Production (MW) values (column B) to the new model output.

$ws.Cells.Item(2, 1).Value2 = 46044.01041666666
$ws.Cells.Item(3, 1).Value2 = 46044.02083333334
$ws.Cells.Item(4, 1).Value2 = 46044.03125
$ws.Cells.Item(5, 1).Value2 = 46044.04166666666
$ws.Cells.Item(6, 1).Value2 = 46044.05208333334
$ws.Cells.Item(6, 2).Value2 = 0.29
$ws.Cells.Item(7, 1).Value2 = 46044.0625
$ws.Cells.Item(8, 1).Value2 = 46044.07291666666
$ws.Cells.Item(9, 1).Value2 = 46044.08333333334
$ws.Cells.Item(10, 1).Value2 = 46044.09375
$ws.Cells.Item(11, 1).Value2 = 46044.10416666666
$ws.Cells.Item(12, 1).Value2 = 46044.11458333334
$ws.Cells.Item(13, 1).Value2 = 46044.125
$ws.Cells.Item(14, 1).Value2 = 46044.13541666666
$ws.Cells.Item(14, 2).Value2 = 0.45
$ws.Cells.Item(15, 1).Value2 = 46044.14583333334
$ws.Cells.Item(16, 1).Value2 = 46044.15625
$ws.Cells.Item(17, 1).Value2 = 46044.16666666666
$ws.Cells.Item(18, 1).Value2 = 46044.17708333334
$ws.Cells.Item(18, 2).Value2 = 0
$ws.Cells.Item(19, 1).Value2 = 46044.1875
$ws.Cells.Item(20, 1).Value2 = 46044.19791666666
$ws.Cells.Item(21, 1).Value2 = 46044.20833333334
$ws.Cells.Item(22, 1).Value2 = 46044.21875
$ws.Cells.Item(22, 2).Value2 = 0.462
$ws.Cells.Item(23, 1).Value2 = 46044.22916666666
$ws.Cells.Item(23, 2).Value2 = 0.47
$ws.Cells.Item(24, 1).Value2 = 46044.23958333334
$ws.Cells.Item(24, 2).Value2 = 0
$ws.Cells.Item(25, 1).Value2 = 46044.25
$ws.Cells.Item(25, 2).Value2 = 0.484
$ws.Cells.Item(26, 1).Value2 = 46044.26041666666
$ws.Cells.Item(26, 2).Value2 = 0.588
$ws.Cells.Item(27, 1).Value2 = 46044.27083333334
$ws.Cells.Item(27, 2).Value2 = 0.841
$ws.Cells.Item(28, 1).Value2 = 46044.28125
$ws.Cells.Item(28, 2).Value2 = 1.197
$ws.Cells.Item(29, 1).Value2 = 46044.29166666666
$ws.Cells.Item(29, 2).Value2 = 2.256
$ws.Cells.Item(30, 1).Value2 = 46044.30208333334
$ws.Cells.Item(30, 2).Value2 = 12.012
$ws.Cells.Item(31, 1).Value2 = 46044.3125
$ws.Cells.Item(31, 2).Value2 = 24.118
$ws.Cells.Item(32, 1).Value2 = 46044.32291666666
$ws.Cells.Item(32, 2).Value2 = 39.222
$ws.Cells.Item(33, 1).Value2 = 46044.33333333334
$ws.Cells.Item(33, 2).Value2 = 56.534
$ws.Cells.Item(34, 1).Value2 = 46044.34375
$ws.Cells.Item(34, 2).Value2 = 93.496
$ws.Cells.Item(35, 1).Value2 = 46044.35416666666
$ws.Cells.Item(35, 2).Value2 = 116.008
$ws.Cells.Item(36, 1).Value2 = 46044.36458333334
$ws.Cells.Item(36, 2).Value2 = 142.186
$ws.Cells.Item(37, 1).Value2 = 46044.375
$ws.Cells.Item(37, 2).Value2 = 166.037
$ws.Cells.Item(38, 1).Value2 = 46044.38541666666
$ws.Cells.Item(38, 2).Value2 = 206.921
$ws.Cells.Item(39, 1).Value2 = 46044.39583333334
$ws.Cells.Item(39, 2).Value2 = 227.295
$ws.Cells.Item(40, 1).Value2 = 46044.40625
$ws.Cells.Item(40, 2).Value2 = 249.076
$ws.Cells.Item(41, 1).Value2 = 46044.41666666666
$ws.Cells.Item(41, 2).Value2 = 265.005
$ws.Cells.Item(42, 1).Value2 = 46044.42708333334
$ws.Cells.Item(42, 2).Value2 = 292.77
$ws.Cells.Item(43, 1).Value2 = 46044.4375
$ws.Cells.Item(43, 2).Value2 = 307.295
$ws.Cells.Item(44, 1).Value2 = 46044.44791666666
$ws.Cells.Item(44, 2).Value2 = 327.611
$ws.Cells.Item(45, 1).Value2 = 46044.45833333334
$ws.Cells.Item(45, 2).Value2 = 337.767
$ws.Cells.Item(46, 1).Value2 = 46044.46875
$ws.Cells.Item(46, 2).Value2 = 353.525
$ws.Cells.Item(47, 1).Value2 = 46044.47916666666
$ws.Cells.Item(47, 2).Value2 = 357.384
$ws.Cells.Item(48, 1).Value2 = 46044.48958333334
$ws.Cells.Item(48, 2).Value2 = 357.203
$ws.Cells.Item(49, 1).Value2 = 46044.5
$ws.Cells.Item(49, 2).Value2 = 356.396
$ws.Cells.Item(50, 1).Value2 = 46044.51041666666
$ws.Cells.Item(50, 2).Value2 = 353.526
$ws.Cells.Item(51, 1).Value2 = 46044.52083333334
$ws.Cells.Item(51, 2).Value2 = 344.782
$ws.Cells.Item(52, 1).Value2 = 46044.53125
$ws.Cells.Item(52, 2).Value2 = 334.905
$ws.Cells.Item(53, 1).Value2 = 46044.54166666666
$ws.Cells.Item(53, 2).Value2 = 320.995
$ws.Cells.Item(54, 1).Value2 = 46044.55208333334
$ws.Cells.Item(54, 2).Value2 = 301.379
$ws.Cells.Item(55, 1).Value2 = 46044.5625
$ws.Cells.Item(55, 2).Value2 = 282.196
$ws.Cells.Item(56, 1).Value2 = 46044.57291666666
$ws.Cells.Item(56, 2).Value2 = 260.077
$ws.Cells.Item(57, 1).Value2 = 46044.58333333334
$ws.Cells.Item(57, 2).Value2 = 237.102
$ws.Cells.Item(58, 1).Value2 = 46044.59375
$ws.Cells.Item(58, 2).Value2 = 191.808
$ws.Cells.Item(59, 1).Value2 = 46044.60416666666
$ws.Cells.Item(59, 2).Value2 = 165.562
$ws.Cells.Item(60, 1).Value2 = 46044.61458333334
$ws.Cells.Item(60, 2).Value2 = 136.962
$ws.Cells.Item(61, 1).Value2 = 46044.625
$ws.Cells.Item(61, 2).Value2 = 112.947
$ws.Cells.Item(62, 1).Value2 = 46044.63541666666
$ws.Cells.Item(62, 2).Value2 = 76.873
$ws.Cells.Item(63, 1).Value2 = 46044.64583333334
$ws.Cells.Item(63, 2).Value2 = 56.2
$ws.Cells.Item(64, 1).Value2 = 46044.65625
$ws.Cells.Item(64, 2).Value2 = 39.025
$ws.Cells.Item(65, 1).Value2 = 46044.66666666666
$ws.Cells.Item(65, 2).Value2 = 26.44
$ws.Cells.Item(66, 1).Value2 = 46044.67708333334
$ws.Cells.Item(66, 2).Value2 = 13.986
$ws.Cells.Item(67, 1).Value2 = 46044.6875
$ws.Cells.Item(67, 2).Value2 = 10.226
$ws.Cells.Item(68, 1).Value2 = 46044.69791666666
$ws.Cells.Item(68, 2).Value2 = 8.977
$ws.Cells.Item(69, 1).Value2 = 46044.70833333334
$ws.Cells.Item(69, 2).Value2 = 8.711
$ws.Cells.Item(70, 1).Value2 = 46044.71875
$ws.Cells.Item(70, 2).Value2 = 2.65
$ws.Cells.Item(71, 1).Value2 = 46044.72916666666
$ws.Cells.Item(72, 1).Value2 = 46044.73958333334
$ws.Cells.Item(73, 1).Value2 = 46044.75
$ws.Cells.Item(73, 2).Value2 = 0
$ws.Cells.Item(74, 1).Value2 = 46044.76041666666
$ws.Cells.Item(74, 2).Value2 = 0
$ws.Cells.Item(75, 1).Value2 = 46044.77083333334
$ws.Cells.Item(75, 2).Value2 = 0.65
$ws.Cells.Item(76, 1).Value2 = 46044.78125
$ws.Cells.Item(76, 2).Value2 = 0
$ws.Cells.Item(77, 1).Value2 = 46044.79166666666
$ws.Cells.Item(78, 1).Value2 = 46044.80208333334
$ws.Cells.Item(78, 2).Value2 = 0.49
$ws.Cells.Item(79, 1).Value2 = 46044.8125
$ws.Cells.Item(80, 1).Value2 = 46044.82291666666
$ws.Cells.Item(81, 1).Value2 = 46044.83333333334
$ws.Cells.Item(82, 1).Value2 = 46044.84375
$ws.Cells.Item(82, 2).Value2 = 0.65
$ws.Cells.Item(83, 1).Value2 = 46044.85416666666
$ws.Cells.Item(84, 1).Value2 = 46044.86458333334
$ws.Cells.Item(85, 1).Value2 = 46044.875
$ws.Cells.Item(86, 1).Value2 = 46044.88541666666
$ws.Cells.Item(86, 2).Value2 = 0.45
$ws.Cells.Item(87, 1).Value2 = 46044.89583333334
$ws.Cells.Item(88, 1).Value2 = 46044.90625
$ws.Cells.Item(89, 1).Value2 = 46044.91666666666
$ws.Cells.Item(90, 1).Value2 = 46044.92708333334
$ws.Cells.Item(91, 1).Value2 = 46044.9375
$ws.Cells.Item(92, 1).Value2 = 46044.94791666666
$ws.Cells.Item(93, 1).Value2 = 46044.95833333334
$ws.Cells.Item(94, 1).Value2 = 46044.96875
$ws.Cells.Item(95, 1).Value2 = 46044.97916666666
$ws.Cells.Item(96, 1).Value2 = 46044.98958333334
$ws.Cells.Item(97, 1).Value2 = 46045
